$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Update the "Date and time" value (B5)
$ws.Range("B5").Value = "2014-08-02 12:53"

# Update the RESULTS table (B10:F22) with the new measured values.
# Row 11: Mean value      -> unchanged
# Row 12: Remove mean
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 0.421637021355784

# Row 13: Differentiate -> unchanged

# Row 14: Square
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 0

# Row 15: Box filter
$ws.Range("C15").Value = 25
$ws.Range("D15").Value = 23
$ws.Range("E15").Value = 30
$ws.Range("F15").Value = 1.82878222991269

# Row 16: Median filter
$ws.Range("C16").Value = 224
$ws.Range("D16").Value = 222
$ws.Range("E16").Value = 231
$ws.Range("F16").Value = 2.64365067451978

# Row 17: Integrate
$ws.Range("C17").Value = 250
$ws.Range("D17").Value = 248
$ws.Range("E17").Value = 256
$ws.Range("F17").Value = 2.60128173535022

# Row 18: Segment
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 0.948683298050514

# Row 19: Treshold
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 12
$ws.Range("F19").Value = 2.20100986922922

# Row 20: Indices
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 16
$ws.Range("E20").Value = 17
$ws.Range("F20").Value = 0.421637021355784

# Row 21: Find peaks -> unchanged

# Row 22: Total time
$ws.Range("C22").Value = 306
$ws.Range("D22").Value = 300
$ws.Range("E22").Value = 309
$ws.Range("F22").Value = 2.85968141193696
